$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.610.73'
$ws.Range("E2").Value = '  -7.08%  '
$ws.Range("D3").Value = '2.433.04'
$ws.Range("E3").Value = '  -10.69%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '468.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.496'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.77%  '
$ws.Range("D9").Value = '2.452.85'
$ws.Range("E9").Value = '  -10.47%  '
$ws.Range("E10").Value = '  -8.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.34'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -12.38%  '
$ws.Range("E12").Value = '  -9.44%  '
$ws.Range("E13").Value = '  -3.89%  '
$ws.Range("D14").Value = '2.854.61'
$ws.Range("E14").Value = '  -10.70%  '
$ws.Range("D15").Value = '54.592.99'
$ws.Range("E15").Value = '  -7.10%  '
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.79%  '
$ws.Range("D18").Value = '2.445.55'
$ws.Range("E18").Value = '  -10.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -13.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").Value = '  -13.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '56.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -10.40%  '
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.158'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.68%  '
$ws.Range("B28").Value = 'Polygon'
$ws.Range("C28").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.388'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.48%  '
$ws.Range("D29").Value = '2.526.73'
$ws.Range("E29").Value = '  -10.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '0.0₃0725'
$ws.Range("E32").Value = '  -13.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("E34").Value = '  -7.53%  '
$ws.Range("E35").Value = '  -10.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.45%  '
$ws.Range("E37").Value = '  -15.12%  '
$ws.Range("E38").Value = '  -6.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.809'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -14.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.994'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '33.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.54%  '
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0528'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.50%  '
$ws.Range("E44").Value = '  -9.09%  '
$ws.Range("E45").Value = '  -10.61%  '
$ws.Range("E46").Value = '  -2.61%  '
$ws.Range("D47").Value = '1.940.89'
$ws.Range("E47").Value = '  -11.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0884'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0219'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '235.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.31%  '
